# Update documentation for recent software enhancements preparing for 11.10.01 release.
# Adds a new "HTML Documentation?" column (E) to the TSTool Command List sheet,
# marks several commands with "Y" in existing / new columns, and refreshes the
# COUNTIF summary row (239).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header (E1): "HTML Documentation?" bold + wrap text ---
$ws.Range("E1").Value = "HTML Documentation?"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").WrapText = $true

# --- Mark existing rows with "Y" in column C / D where newly required ---
$ws.Range("D71").Value = "Y"

$ws.Range("C81").Value = "Y"

$ws.Range("C118").Value = "Y"
$ws.Range("D118").Value = "Y"
$ws.Range("E118").Value = "Y"

$ws.Range("C128").Value = "Y"
$ws.Range("D128").Value = "Y"

$ws.Range("E186").Value = "Y"

$ws.Range("E213").Value = "Y"

$ws.Range("C236").Value = "Y"
$ws.Range("D236").Value = "Y"

# --- New COUNTIF summary formula for column E (row 239) ---
$ws.Range("E239").Formula = '=COUNTIF(E2:E238,"=Y")+COUNTIF(E2:E238,"=NA")+COUNTIF(E2:E238,"=Z")'
$ws.Range("E239").HorizontalAlignment = -4108

# --- Update the view: active selection / frozen-pane state ---
$ws.Range("E239").Select()
